$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (character units; closest achievable match to target 32 / 55.28515625 / 67.28515625)
$ws.Columns.Item(1).ColumnWidth = 31.16
$ws.Columns.Item(2).ColumnWidth = 54.5
$ws.Columns.Item(3).ColumnWidth = 66.5

# Row 1
$ws.Range("A1").Value = "cd"
$ws.Range("B1").Value = "cd OR cd any-path"
$ws.Range("C1").Value = "للوصول الى المجلد root للجهاز او للوصول الى مسار محدد"

# Row 2
$ws.Range("A2").Value = "mkdir"
$ws.Range("B2").Value = "mkdir name-of-folder"
$ws.Range("C2").Value = "لإنشاء مجلد في المسار الحالي"

# Row 3
$ws.Range("A3").Value = "git clone"
$ws.Range("B3").Value = "git clone any-path"
$ws.Range("C3").Value = "نسخ المستودع الذي تم تحديدة في المسار الى المسار الحالي"

# Row 4
$ws.Range("A4").Value = "dir"
$ws.Range("C4").Value = "يقوم بأستعراض الملفات الموجودة في المسار الحالي"

# Row 5
$ws.Range("A5").Value = "git status"
$ws.Range("C5").Value = "يوضح حالة working dirctory وماذا يوجد بها من تغييرات"

# Row 6
$ws.Range("A6").Value = "git add"
$ws.Range("B6").Value = "git add * OR git add name_file _in_working_dirctory"
$ws.Range("C6").Value = "نقل الملفات من working dirctory الى staging area"

# Row 7
$ws.Range("A7").Value = "git reset head "
$ws.Range("B7").Value = "[git reset head name_file] OR [git restore --staged name_file]"
$ws.Range("C7").Value = "التراجع عن الملفات التي في منطقة staging area وارجاعها الى منطقة working dirctory"

# Row 8
$ws.Range("A8").Value = "git commit -m `"msg-text`""

# Selection matches target xml (activeCell A9)
$ws.Range("A9").Select()
